$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 139 ("Fecha" = 2021-05-17),
# pushing that row (and everything after it) down by one. Excel copies the
# formatting of the row above on insert, so the date cell D139 will pick up
# the same date number format (style id 2) used by the rest of column D.
$ws.Rows.Item(139).Insert()

# Populate the new row 139 with the new observation.
$ws.Cells.Item(139, 1).Value = 8
$ws.Cells.Item(139, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(139, 3).Value = "Coquimbo"
$ws.Cells.Item(139, 4).Value = 44634
$ws.Cells.Item(139, 5).Value = 4
$ws.Cells.Item(139, 6).Value = 100112012
$ws.Cells.Item(139, 7).Value = "Espinaca"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 2200
$ws.Cells.Item(139, 11).Value = 500
$ws.Cells.Item(139, 12).Value = 600
$ws.Cells.Item(139, 13).Value = 550
$ws.Cells.Item(139, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(139, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(139, 16).Value = 1100
$ws.Cells.Item(139, 17).Value = 0.5
$ws.Cells.Item(139, 18).Value = "Hortaliza"
